$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D ("Branch NicName"), shifting AccountName/AccountType/Currency right.
$ws.Columns("D:D").Insert()

# The inserted column has no explicit width yet; set it as close as the engine allows to 10.75
# (mirrors column C's width in the target file).
$ws.Columns("D").ColumnWidth = 9.75

# Header row (re-stamp every header so the table's column-name cache, which
# is keyed by position and not by live cell content, stays in sync with the
# shifted columns).
$ws.Range("D1").Value = "Branch NicName"
$ws.Range("E1").Value = "AccountName"
$ws.Range("F1").Value = "AccountType"
$ws.Range("G1").Value = "Currency"

# Data rows: nickname derived from the Branch column (now in column C)
$ws.Range("D2").Value = "(SG)"
$ws.Range("D3").Value = "(SG)"
$ws.Range("D4").Value = "(BK)"
$ws.Range("D5").Value = "(SG)"
$ws.Range("D6").Value = "(SG)"
$ws.Range("D7").Value = "(SG)"
$ws.Range("D8").Value = "(SG)"
$ws.Range("D9").Value = "(BK)"

# Grow the table (Table1) to include the new column.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:G9"))

# Re-set the last header cell once more after the resize so the newly
# created 7th ListColumn (auto-named "Column7") picks up "Currency".
$ws.Range("G1").Value = "Currency"
